# Convert the old "bookmarked Heading1 + bold byline" header into a
# pandoc-style title block: a Title-styled paragraph (one run per word /
# separator, matching pandoc's output) followed by an Authors-styled
# paragraph, with the wrapping bookmark removed.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the old "Winter Arrives, Work in House And Trip to West"
#    Heading1 paragraph completely (text + paragraph mark). Doing this
#    on the paragraph's own Range deletes it cleanly and, as a side
#    effect, collapses the bookmarkStart/bookmarkEnd pair that used to
#    straddle it down to the same (now empty) position at the very
#    start of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Delete()

# ------------------------------------------------------------------
# 2. The bookmarkStart/bookmarkEnd markers are now both sitting,
#    adjacent and zero-width, at document position 0. A collapsed
#    Range.Delete() right there removes one marker at a time without
#    touching any real text, so call it twice to clear both.
# ------------------------------------------------------------------
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# ------------------------------------------------------------------
# 3. Insert the new Title paragraph (one run per word/space/punct,
#    mirroring how pandoc emits its title block) at the very start.
# ------------------------------------------------------------------
$titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
            "<w:pPr><w:pStyle w:val='Title'/></w:pPr>" +
            "<w:r><w:t xml:space='preserve'>Winter</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>Arrives</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>,</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>Work</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>in</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>House</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>And</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>Trip</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>to</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>West</w:t></w:r>" +
            "</w:p>"
$d.Range(0, 0).InsertXML($titleXml)

# ------------------------------------------------------------------
# 4. Replace the old bold "By Dorothy Day" paragraph (now paragraph 2)
#    with an Authors-styled paragraph split into "Dorothy" / " " /
#    "Day" runs, dropping the leading "By ".
# ------------------------------------------------------------------
$authorPara = $d.Paragraphs.Item(2)
$authorXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
             "<w:pPr><w:pStyle w:val='Authors'/></w:pPr>" +
             "<w:r><w:t xml:space='preserve'>Dorothy</w:t></w:r>" +
             "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
             "<w:r><w:t xml:space='preserve'>Day</w:t></w:r>" +
             "</w:p>"
$authorPara.Range.InsertXML($authorXml)

Write-Output ("paragraph 1: [" + $d.Paragraphs.Item(1).Range.Text + "] style=" + $d.Paragraphs.Item(1).Range.Style.NameLocal)
Write-Output ("paragraph 2: [" + $d.Paragraphs.Item(2).Range.Text + "] style=" + $d.Paragraphs.Item(2).Range.Style.NameLocal)
